$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow a few columns (J=10, AA=27, AB=28) from width 8 to width 7 ---
# ColumnWidth uses character units; stored OOXML width differs by a small
# fixed offset. 6.15 lands safely inside the pixel-rounding band that
# serializes back out as width="7".
$ws.Columns.Item(10).ColumnWidth = 6.15
$ws.Columns.Item(27).ColumnWidth = 6.15
$ws.Columns.Item(28).ColumnWidth = 6.15

# --- Reduce row 5's numeric precision to 2 decimal places (custom accuracy) ---
$ws.Range("B5").Value = 18.64
$ws.Range("C5").Value = 13.83
$ws.Range("D5").Value = 1.17
$ws.Range("E5").Value = 40.76
$ws.Range("F5").Value = 33.06
$ws.Range("G5").Value = 14.61
$ws.Range("H5").Value = 56.49
$ws.Range("I5").Value = 22.69
$ws.Range("J5").Value = 10.03
$ws.Range("K5").Value = 14.78
$ws.Range("L5").Value = 16.33
$ws.Range("M5").Value = 17.38
$ws.Range("N5").Value = 4.71
$ws.Range("O5").Value = 14.66
$ws.Range("P5").Value = 20.82
$ws.Range("Q5").Value = 12.44
$ws.Range("R5").Value = 0.78
$ws.Range("S5").Value = 0.79
$ws.Range("T5").Value = 215.86
$ws.Range("U5").Value = 41
$ws.Range("V5").Value = 13.53
$ws.Range("W5").Value = 27.48
$ws.Range("X5").Value = 14.39
$ws.Range("Y5").Value = 2.23
$ws.Range("Z5").Value = 27.63
$ws.Range("AA5").Value = 11.95
$ws.Range("AB5").Value = 10.63
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 17.12
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 51.45
$ws.Range("AG5").Value = 7.58
$ws.Range("AH5").Value = 16.92

# --- Drop row 6 entirely (data trimmed from 1000 rows) ---
$ws.Rows.Item(6).Delete()
